# -----------------------------------------------------------------------
# Fatture.xlsx fix-up ("had a bug, fixed it"):
#   1. Remove the stray empty "Sconto o magg" (H) placeholder cells that
#      were left on the DIAGRES "DOD *CASSANI" lines in both the
#      "Acquista" and "Inventario" sheets.
#   2. Append the invoice rows that were missing from the bottom of both
#      sheets (ADRIATICA LUCIDANTI DI PETRIZZO GIANLUCA C SNC, NP DIAMANT
#      SRLS, and SWEDIAM SRL).
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$wsAcquista = $wb.Worksheets.Item("Acquista")
$wsInventario = $wb.Worksheets.Item("Inventario")

# --- Remove the empty inline-string placeholder cells in column H ---
# "Acquista" sheet: rows 323-334
foreach ($r in 323..334) {
    $wsAcquista.Range("H$r").ClearContents()
}

# "Inventario" sheet: rows 322-333 (same data, shifted up by one row)
foreach ($r in 322..333) {
    $wsInventario.Range("H$r").ClearContents()
}

# --- Append the missing rows to "Acquista" ---
$wsAcquista.Range("A352").Value = "ADRIATICA`tLUCIDANTI`tDI`tPETRIZZO`nGIANLUCA`tC`tSNC`n"
$wsAcquista.Range("B352").Value = "31-03-2021"
$wsAcquista.Rows.Item(352).EntireRow.AutoFit()
$wsAcquista.Range("D354").Value = "LUCIDANTE`tPOTE`tGRANITO"
$wsAcquista.Range("E354").Value = "5,00"
$wsAcquista.Range("F354").Value = "13,60"
$wsAcquista.Range("G354").Value = "KG"
$wsAcquista.Range("I354").Value = "22,00"
$wsAcquista.Range("J354").Value = "68,00"
$wsAcquista.Range("D355").Value = "CASSANI`tSUPER`tGR`t30/60/120"
$wsAcquista.Range("E355").Value = "135,00"
$wsAcquista.Range("F355").Value = "0,97"
$wsAcquista.Range("G355").Value = "PZ"
$wsAcquista.Range("I355").Value = "22,00"
$wsAcquista.Range("J355").Value = "130,95"
$wsAcquista.Range("D356").Value = "LUCIDANTE`tF1"
$wsAcquista.Range("E356").Value = "50,00"
$wsAcquista.Range("F356").Value = "4,30"
$wsAcquista.Range("G356").Value = "KG"
$wsAcquista.Range("I356").Value = "22,00"
$wsAcquista.Range("J356").Value = "215,00"
$wsAcquista.Range("D357").Value = "LUCIDANTE`tWR`tMARMO"
$wsAcquista.Range("E357").Value = "50,00"
$wsAcquista.Range("F357").Value = "4,60"
$wsAcquista.Range("G357").Value = "KG"
$wsAcquista.Range("I357").Value = "22,00"
$wsAcquista.Range("J357").Value = "230,00"
$wsAcquista.Range("A359").Value = "NP`tDIAMANT`tSRLS`n"
$wsAcquista.Range("B359").Value = "27-02-2021"
$wsAcquista.Rows.Item(359).EntireRow.AutoFit()
$wsAcquista.Range("C361").Value = "X"
$wsAcquista.Range("D361").Value = "DISCO`tD.350`tB2DLC/2`tSILENZIATO`tGRANITO`tF.50"
$wsAcquista.Range("E361").Value = "1,00"
$wsAcquista.Range("F361").Value = "51,00"
$wsAcquista.Range("G361").Value = "NR"
$wsAcquista.Range("I361").Value = "22,00"
$wsAcquista.Range("J361").Value = "51,00"
$wsAcquista.Range("C362").Value = "X"
$wsAcquista.Range("D362").Value = "DISCO`tD.400`tB2DLC/2`tSILENZIATO`tGRANITO"
$wsAcquista.Range("E362").Value = "2,00"
$wsAcquista.Range("F362").Value = "67,00"
$wsAcquista.Range("G362").Value = "NR"
$wsAcquista.Range("I362").Value = "22,00"
$wsAcquista.Range("J362").Value = "134,00"
$wsAcquista.Range("C363").Value = "X"
$wsAcquista.Range("D363").Value = "DISCO`tD.500`tB2DLC/2`tSILENZIATO`tGRANITO"
$wsAcquista.Range("E363").Value = "1,00"
$wsAcquista.Range("F363").Value = "106,00"
$wsAcquista.Range("G363").Value = "NR"
$wsAcquista.Range("I363").Value = "22,00"
$wsAcquista.Range("J363").Value = "106,00"
$wsAcquista.Range("C364").Value = "X"
$wsAcquista.Range("D364").Value = "DISCO`tD.500`tM1S`tSILENZIATO`tMARMO`tF.60/50"
$wsAcquista.Range("E364").Value = "1,00"
$wsAcquista.Range("F364").Value = "112,00"
$wsAcquista.Range("G364").Value = "NR"
$wsAcquista.Range("I364").Value = "22,00"
$wsAcquista.Range("J364").Value = "112,00"
$wsAcquista.Range("C365").Value = "X"
$wsAcquista.Range("D365").Value = "DISCO`tD.625`tV55DLC/2`tSILENZIATO`tPER"
$wsAcquista.Range("E365").Value = "1,00"
$wsAcquista.Range("F365").Value = "245,10"
$wsAcquista.Range("G365").Value = "NR"
$wsAcquista.Range("I365").Value = "22,00"
$wsAcquista.Range("J365").Value = "245,10"
$wsAcquista.Range("C366").Value = "X"
$wsAcquista.Range("D366").Value = "FRESA`tDA`tTAGLIO`t22*45`tATT.1/2`tGAS`t5`tSETTORI"
$wsAcquista.Range("E366").Value = "2,00"
$wsAcquista.Range("F366").Value = "50,00"
$wsAcquista.Range("G366").Value = "NR"
$wsAcquista.Range("I366").Value = "22,00"
$wsAcquista.Range("J366").Value = "100,00"
$wsAcquista.Range("C367").Value = "X"
$wsAcquista.Range("D367").Value = "DISCO`tD.230`tD.E.`tF.25,4`tMARMO"
$wsAcquista.Range("E367").Value = "1,00"
$wsAcquista.Range("F367").Value = "28,00"
$wsAcquista.Range("G367").Value = "NR"
$wsAcquista.Range("I367").Value = "22,00"
$wsAcquista.Range("J367").Value = "28,00"
$wsAcquista.Range("C368").Value = "X"
$wsAcquista.Range("D368").Value = "MOLA`tD.240X437`tRIGENERATA`tCONH`t220"
$wsAcquista.Range("E368").Value = "1,00"
$wsAcquista.Range("F368").Value = "619,00"
$wsAcquista.Range("G368").Value = "NR"
$wsAcquista.Range("I368").Value = "22,00"
$wsAcquista.Range("J368").Value = "619,00"
$wsAcquista.Range("A370").Value = "SWEDIAM`tSRL`n"
$wsAcquista.Range("B370").Value = "29-01-2021"
$wsAcquista.Rows.Item(370).EntireRow.AutoFit()
$wsAcquista.Range("A372").Value = "SWEDIAM`tSRL`n"
$wsAcquista.Range("B372").Value = "29-01-2021"
$wsAcquista.Rows.Item(372).EntireRow.AutoFit()
$wsAcquista.Range("A374").Value = "SWEDIAM`tSRL`n"
$wsAcquista.Range("B374").Value = "26-02-2021"
$wsAcquista.Rows.Item(374).EntireRow.AutoFit()

# --- Append the missing rows to "Inventario" ---
$wsInventario.Range("D353").Value = "LUCIDANTE`tPOTE`tGRANITO"
$wsInventario.Range("E353").Value = "5,00"
$wsInventario.Range("F353").Value = "13,60"
$wsInventario.Range("G353").Value = "KG"
$wsInventario.Range("I353").Value = "22,00"
$wsInventario.Range("J353").Value = "68,00"
$wsInventario.Range("D354").Value = "CASSANI`tSUPER`tGR`t30/60/120"
$wsInventario.Range("E354").Value = "135,00"
$wsInventario.Range("F354").Value = "0,97"
$wsInventario.Range("G354").Value = "PZ"
$wsInventario.Range("I354").Value = "22,00"
$wsInventario.Range("J354").Value = "130,95"
$wsInventario.Range("D355").Value = "LUCIDANTE`tF1"
$wsInventario.Range("E355").Value = "50,00"
$wsInventario.Range("F355").Value = "4,30"
$wsInventario.Range("G355").Value = "KG"
$wsInventario.Range("I355").Value = "22,00"
$wsInventario.Range("J355").Value = "215,00"
$wsInventario.Range("D356").Value = "LUCIDANTE`tWR`tMARMO"
$wsInventario.Range("E356").Value = "50,00"
$wsInventario.Range("F356").Value = "4,60"
$wsInventario.Range("G356").Value = "KG"
$wsInventario.Range("I356").Value = "22,00"
$wsInventario.Range("J356").Value = "230,00"
$wsInventario.Range("C360").Value = "X"
$wsInventario.Range("D360").Value = "DISCO`tD.350`tB2DLC/2`tSILENZIATO`tGRANITO`tF.50"
$wsInventario.Range("E360").Value = "1,00"
$wsInventario.Range("F360").Value = "51,00"
$wsInventario.Range("G360").Value = "NR"
$wsInventario.Range("I360").Value = "22,00"
$wsInventario.Range("J360").Value = "51,00"
$wsInventario.Range("C361").Value = "X"
$wsInventario.Range("D361").Value = "DISCO`tD.400`tB2DLC/2`tSILENZIATO`tGRANITO"
$wsInventario.Range("E361").Value = "2,00"
$wsInventario.Range("F361").Value = "67,00"
$wsInventario.Range("G361").Value = "NR"
$wsInventario.Range("I361").Value = "22,00"
$wsInventario.Range("J361").Value = "134,00"
$wsInventario.Range("C362").Value = "X"
$wsInventario.Range("D362").Value = "DISCO`tD.500`tB2DLC/2`tSILENZIATO`tGRANITO"
$wsInventario.Range("E362").Value = "1,00"
$wsInventario.Range("F362").Value = "106,00"
$wsInventario.Range("G362").Value = "NR"
$wsInventario.Range("I362").Value = "22,00"
$wsInventario.Range("J362").Value = "106,00"
$wsInventario.Range("C363").Value = "X"
$wsInventario.Range("D363").Value = "DISCO`tD.500`tM1S`tSILENZIATO`tMARMO`tF.60/50"
$wsInventario.Range("E363").Value = "1,00"
$wsInventario.Range("F363").Value = "112,00"
$wsInventario.Range("G363").Value = "NR"
$wsInventario.Range("I363").Value = "22,00"
$wsInventario.Range("J363").Value = "112,00"
$wsInventario.Range("C364").Value = "X"
$wsInventario.Range("D364").Value = "DISCO`tD.625`tV55DLC/2`tSILENZIATO`tPER"
$wsInventario.Range("E364").Value = "1,00"
$wsInventario.Range("F364").Value = "245,10"
$wsInventario.Range("G364").Value = "NR"
$wsInventario.Range("I364").Value = "22,00"
$wsInventario.Range("J364").Value = "245,10"
$wsInventario.Range("C365").Value = "X"
$wsInventario.Range("D365").Value = "FRESA`tDA`tTAGLIO`t22*45`tATT.1/2`tGAS`t5`tSETTORI"
$wsInventario.Range("E365").Value = "2,00"
$wsInventario.Range("F365").Value = "50,00"
$wsInventario.Range("G365").Value = "NR"
$wsInventario.Range("I365").Value = "22,00"
$wsInventario.Range("J365").Value = "100,00"
$wsInventario.Range("C366").Value = "X"
$wsInventario.Range("D366").Value = "DISCO`tD.230`tD.E.`tF.25,4`tMARMO"
$wsInventario.Range("E366").Value = "1,00"
$wsInventario.Range("F366").Value = "28,00"
$wsInventario.Range("G366").Value = "NR"
$wsInventario.Range("I366").Value = "22,00"
$wsInventario.Range("J366").Value = "28,00"
$wsInventario.Range("C367").Value = "X"
$wsInventario.Range("D367").Value = "MOLA`tD.240X437`tRIGENERATA`tCONH`t220"
$wsInventario.Range("E367").Value = "1,00"
$wsInventario.Range("F367").Value = "619,00"
$wsInventario.Range("G367").Value = "NR"
$wsInventario.Range("I367").Value = "22,00"
$wsInventario.Range("J367").Value = "619,00"
